$d = $word.ActiveDocument

$pairs = @(
    @("494×2=", "277×7="),
    @("800×9=", "855×8="),
    @("267×3=", "812×4="),
    @("186×8=", "347×7="),
    @("664×4=", "964×5="),
    @("721×9=", "932×8="),
    @("770×2=", "144×8="),
    @("199×4=", "716×5="),
    @("554×8=", "235×8="),
    @("542×2=", "564×7="),
    @("616×4=", "355×8="),
    @("442×4=", "449×9="),
    @("590×4=", "462×9="),
    @("523×7=", "457×3="),
    @("960×2=", "256×9="),
    @("393×2=", "342×8="),
    @("693×8=", "963×5="),
    @("714×5=", "389×9="),
    @("582×6=", "691×3="),
    @("156×8=", "609×8="),
    @("950×5=", "334×8="),
    @("950×4=", "702×3="),
    @("190×8=", "139×6="),
    @("849×7=", "445×7="),
    @("147×8=", "958×3=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
